# Auto-generated: applies cell value updates per the commit diff
# (cryptos.xlsx price/volume refresh, GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.611.81'
$ws.Range('E2').Value = '  +3.38%  '
$ws.Range('D3').Value = '1.699.32'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3945'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4045'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('B9').Value = 'Polygon'
$ws.Range('C9').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.555'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.59%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +15.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.003'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08812'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.279'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +11.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001333'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.700'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.26%  '
$ws.Range('D17').Value = '1.698.90'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '101.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07042'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.926'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').Value = '24.558.16'
$ws.Range('E24').Value = '  +3.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.976'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.341'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.248'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.608'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +28.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.115'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = '1.886.69'
$ws.Range('E33').Value = '  +3.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.578'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08550'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.987'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2748'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02786'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09056'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.470'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7740'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7289'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.507'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.189'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.294'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.25%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08007'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.35%  '
